$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3089.6875
$ws.Range("I64").Value = 3017.55
$ws.Range("J64").Value = 3450.375
$ws.Range("K64").Value = 3017.55
$ws.Range("L64").Value = 3450.375
$ws.Range("M64").Value = -2769.55
$ws.Range("N64").Value = -3946.375
$ws.Range("H67").Value = 3089.6875
$ws.Range("I67").Value = 3017.55
$ws.Range("J67").Value = 3450.375
$ws.Range("K67").Value = 3017.55
$ws.Range("L67").Value = 3450.375
$ws.Range("M67").Value = -2159.55
$ws.Range("N67").Value = -5166.375
$ws.Range("H129").Value = 1158.5625
$ws.Range("I129").Value = 356.2857
$ws.Range("J129").Value = 1383.2
$ws.Range("K129").Value = 1068.8571
$ws.Range("L129").Value = 4149.6
$ws.Range("M129").Value = 3931.1429
$ws.Range("N129").Value = -14149.6

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20758.629
$ws.Range("I32").Value = 4131.7617
$ws.Range("J32").Value = 90591.47
$ws.Range("K32").Value = 4131.7617
$ws.Range("L32").Value = 90591.47
$ws.Range("M32").Value = -3844.7617
$ws.Range("N32").Value = -91165.47
$ws.Range("H80").Value = 30805
$ws.Range("J80").Value = 30805
$ws.Range("L80").Value = 30805
$ws.Range("N80").Value = -32801
$ws.Range("H83").Value = 30805
$ws.Range("J83").Value = 30805
$ws.Range("L83").Value = 92415
$ws.Range("N83").Value = -102399
$ws.Range("H132").Value = 2069.5557
$ws.Range("I132").Value = 1855.0322
$ws.Range("K132").Value = 5565.096600000001
$ws.Range("M132").Value = -3035.096600000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 23400
$ws.Range("J76").Value = 23400
$ws.Range("L76").Value = 23400
$ws.Range("N76").Value = -24030
$ws.Range("H79").Value = 23400
$ws.Range("J79").Value = 23400
$ws.Range("L79").Value = 23400
$ws.Range("N79").Value = -25584
$ws.Range("H105").Value = 2959.7646
$ws.Range("I105").Value = 1892.0834
$ws.Range("J105").Value = 5522.2
$ws.Range("K105").Value = 1892.0834
$ws.Range("L105").Value = 5522.2
$ws.Range("M105").Value = -145.0834
$ws.Range("N105").Value = -9016.200000000001
$ws.Range("H134").Value = 2260.6667
$ws.Range("I134").Value = 2152.5925
$ws.Range("J134").Value = 3233.3333
$ws.Range("K134").Value = 6457.7775
$ws.Range("L134").Value = 9699.999899999999
$ws.Range("M134").Value = -3922.7775
$ws.Range("N134").Value = -14769.9999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15046
$ws.Range("J50").Value = 15046
$ws.Range("L50").Value = 15046
$ws.Range("N50").Value = -16296
$ws.Range("H51").Value = 19891.777
$ws.Range("J51").Value = 19891.777
$ws.Range("L51").Value = 19891.777
$ws.Range("N51").Value = -21363.777
$ws.Range("H61").Value = 19891.777
$ws.Range("J61").Value = 19891.777
$ws.Range("L61").Value = 19891.777
$ws.Range("N61").Value = -20587.777
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H74").Value = 32545
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 32545
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 32545
$ws.Range("N74").Value = -34293
$ws.Range("H77").Value = 32545
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 32545
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 97635
$ws.Range("N77").Value = -106371
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 915.41174
$ws.Range("I107").Value = 572.2778
$ws.Range("J107").Value = 1738.9333
$ws.Range("K107").Value = 1716.8334
$ws.Range("L107").Value = 5216.7999
$ws.Range("M107").Value = 203.1666
$ws.Range("N107").Value = -9056.7999
$ws.Range("H131").Value = 1026.5358
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1026.5358
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 3079.6074
$ws.Range("N131").Value = -13159.6074
$ws.Range("M131").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4326.933
$ws.Range("I122").Value = 3531.077
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 10593.231
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -8143.231
$ws.Range("N122").Value = -33400
$ws.Range("H126").Value = 4685.4146
$ws.Range("I126").Value = 2416.1667
$ws.Range("K126").Value = 7248.500100000001
$ws.Range("M126").Value = -4778.500100000001
$ws.Range("H140").Value = 59000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 59000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 59000
$ws.Range("N140").Value = -69360
$ws.Range("M140").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3633.1667
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3633.1667
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3633.1667
$ws.Range("N7").Value = -3857.1667
$ws.Range("H22").Value = 812.8421
$ws.Range("I22").Value = 674.75
$ws.Range("J22").Value = 913.2727
$ws.Range("K22").Value = 674.75
$ws.Range("L22").Value = 913.2727
$ws.Range("M22").Value = -379.75
$ws.Range("N22").Value = -1503.2727
$ws.Range("H27").Value = 812.8421
$ws.Range("I27").Value = 674.75
$ws.Range("J27").Value = 913.2727
$ws.Range("K27").Value = 674.75
$ws.Range("L27").Value = 913.2727
$ws.Range("M27").Value = -567.75
$ws.Range("N27").Value = -1127.2727
$ws.Range("H40").Value = 4997.25
$ws.Range("I40").Value = 4996.3335
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4996.3335
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4860.3335
$ws.Range("N40").Value = -5272
$ws.Range("H126").Value = 3633.1667
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3633.1667
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10899.5001
$ws.Range("N126").Value = -15839.5001
$ws.Range("M7").ClearContents()
$ws.Range("M126").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 13750
$ws.Range("J21").Value = 13750
$ws.Range("L21").Value = 13750
$ws.Range("N21").Value = -14220
$ws.Range("H35").Value = 13750
$ws.Range("J35").Value = 13750
$ws.Range("L35").Value = 13750
$ws.Range("N35").Value = -14330
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H82").Value = 38375
$ws.Range("J82").Value = 38375
$ws.Range("L82").Value = 38375
$ws.Range("N82").Value = -39141
$ws.Range("H85").Value = 38375
$ws.Range("J85").Value = 38375
$ws.Range("L85").Value = 38375
$ws.Range("N85").Value = -41027
$ws.Range("H100").Value = 838087.3
$ws.Range("I100").Value = 6250
$ws.Range("J100").Value = 3333599.2
$ws.Range("K100").Value = 12500
$ws.Range("L100").Value = 6667198.4
$ws.Range("M100").Value = -11959
$ws.Range("N100").Value = -6668280.4
$ws.Range("H107").Value = 2851.6667
$ws.Range("I107").Value = 2366.6428
$ws.Range("J107").Value = 4549.25
$ws.Range("K107").Value = 7099.928400000001
$ws.Range("L107").Value = 13647.75
$ws.Range("M107").Value = -5179.928400000001
$ws.Range("N107").Value = -17487.75
$ws.Range("H109").Value = 37800
$ws.Range("J109").Value = 37800
$ws.Range("L109").Value = 37800
$ws.Range("N109").Value = -40574
$ws.Range("H122").Value = 3504.5217
$ws.Range("I122").Value = 3413.4
$ws.Range("J122").Value = 3675.375
$ws.Range("K122").Value = 10240.2
$ws.Range("L122").Value = 11026.125
$ws.Range("M122").Value = -7790.200000000001
$ws.Range("N122").Value = -15926.125
$ws.Range("H126").Value = 2556.6428
$ws.Range("I126").Value = 1837.875
$ws.Range("J126").Value = 3515
$ws.Range("K126").Value = 5513.625
$ws.Range("L126").Value = 10545
$ws.Range("M126").Value = -3043.625
$ws.Range("N126").Value = -15485
$ws.Range("H132").Value = 2492.5425
$ws.Range("I132").Value = 2652.35
$ws.Range("K132").Value = 7957.049999999999
$ws.Range("M132").Value = -5427.049999999999
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()
